$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 28.72393179789144
$ws.Range("C2").Value = 9.747167794820477
$ws.Range("D2").Value = 4.276585547397751
$ws.Range("E2").Value = 9.78211633705458
$ws.Range("F2").Value = 66.8069703912669
$ws.Range("J2").Value = 10.43950258765046
$ws.Range("L2").Value = 12.05500587233376

$ws.Range("B3").Value = 28.61399692459365
$ws.Range("C3").Value = 9.534477755564938
$ws.Range("D3").Value = 4.146541706474164
$ws.Range("E3").Value = 9.794217008429278
$ws.Range("F3").Value = 66.08808740486738
$ws.Range("J3").Value = 10.44178200816692
$ws.Range("L3").Value = 12.0963943050855

$ws.Range("B4").Value = 28.55689149624337
$ws.Range("C4").Value = 9.4070637022956
$ws.Range("D4").Value = 4.064582800197979
$ws.Range("E4").Value = 9.802104641717648
$ws.Range("F4").Value = 65.6518923055331
$ws.Range("J4").Value = 10.44364368826155
$ws.Range("L4").Value = 12.1242638691282

$ws.Range("B5").Value = 28.5362516954
$ws.Range("C5").Value = 9.35603993239258
$ws.Range("D5").Value = 4.030685727209845
$ws.Range("E5").Value = 9.805434356480102
$ws.Range("F5").Value = 65.47556705947579
$ws.Range("J5").Value = 10.44451826866731
$ws.Range("L5").Value = 12.13623850365967

$ws.Range("B6").Value = 28.53298380999414
$ws.Range("C6").Value = 9.347624623434696
$ws.Range("D6").Value = 4.025028065369259
$ws.Range("E6").Value = 9.805994235110377
$ws.Range("F6").Value = 65.44637804293765
$ws.Range("J6").Value = 10.4446704860639
$ws.Range("L6").Value = 12.13826417438754

$ws.Range("B7").Value = 28.5566024674948
$ws.Range("C7").Value = 9.406371809385734
$ws.Range("D7").Value = 4.064127624220835
$ws.Range("E7").Value = 9.802149079559548
$ws.Range("F7").Value = 65.64950838521575
$ws.Range("J7").Value = 10.4436550141152
$ws.Range("L7").Value = 12.12442286306496

$ws.Range("B8").Value = 28.68388029844533
$ws.Range("C8").Value = 9.673232263808956
$ws.Range("D8").Value = 4.232203536720889
$ws.Range("E8").Value = 9.786193835927357
$ws.Range("F8").Value = 66.55809279230863
$ws.Range("J8").Value = 10.44019245891575
$ws.Range("L8").Value = 12.06876652887398

$ws.Range("B9").Value = 29.01500827744448
$ws.Range("C9").Value = 10.21740486231156
$ws.Range("D9").Value = 4.543739281661292
$ws.Range("E9").Value = 9.758523066601517
$ws.Range("F9").Value = 68.37569425621066
$ws.Range("J9").Value = 10.43708166561323
$ws.Range("L9").Value = 11.97913038574788

$ws.Range("B10").Value = 29.30647077796345
$ws.Range("C10").Value = 10.62419420524055
$ws.Range("D10").Value = 4.760050035971878
$ws.Range("E10").Value = 9.740378253038157
$ws.Range("F10").Value = 69.72556874423762
$ws.Range("J10").Value = 10.43705738439154
$ws.Range("L10").Value = 11.92518198382263

$ws.Range("B11").Value = 29.4491266366867
$ws.Range("C11").Value = 10.80961680102557
$ws.Range("D11").Value = 4.855444979788985
$ws.Range("E11").Value = 9.73259385131208
$ws.Range("F11").Value = 70.34120785380061
$ws.Range("J11").Value = 10.43754127291924
$ws.Range("L11").Value = 11.9032283706032

$ws.Range("B12").Value = 29.50455575985173
$ws.Range("C12").Value = 10.87979501633361
$ws.Range("D12").Value = 4.891116087258959
$ws.Range("E12").Value = 9.729713325611938
$ws.Range("F12").Value = 70.57442545900138
$ws.Range("J12").Value = 10.43779597958884
$ws.Range("L12").Value = 11.89528750099576

$ws.Range("B13").Value = 29.49255612304242
$ws.Range("C13").Value = 10.86468393055581
$ws.Range("D13").Value = 4.883454123037753
$ws.Range("E13").Value = 9.730330711662329
$ws.Range("F13").Value = 70.52419612836013
$ws.Range("J13").Value = 10.43773794065435
$ws.Range("L13").Value = 11.8969811371718

$ws.Range("B14").Value = 29.45365883250081
$ws.Range("C14").Value = 10.81539165394356
$ws.Range("D14").Value = 4.858388849999115
$ws.Range("E14").Value = 9.732355522438548
$ws.Range("F14").Value = 70.36039382115842
$ws.Range("J14").Value = 10.43756079375472
$ws.Range("L14").Value = 11.90256760323442

$ws.Range("B15").Value = 29.43001532191312
$ws.Range("C15").Value = 10.78519124852997
$ws.Range("D15").Value = 4.842976093273665
$ws.Range("E15").Value = 9.733604527880747
$ws.Range("F15").Value = 70.26006761568
$ws.Range("J15").Value = 10.43746160228429
$ws.Range("L15").Value = 11.90603799447758

$ws.Range("B16").Value = 29.29734718437932
$ws.Range("C16").Value = 10.61207600965674
$ws.Range("D16").Value = 4.753753448949278
$ws.Range("E16").Value = 9.74089640970641
$ws.Range("F16").Value = 69.68535515068059
$ws.Range("J16").Value = 10.43703574597403
$ws.Range("L16").Value = 11.92666880861482

$ws.Range("B17").Value = 29.21851124231738
$ws.Range("C17").Value = 10.50590448819058
$ws.Range("D17").Value = 4.698232947265939
$ws.Range("E17").Value = 9.745489851147905
$ws.Range("F17").Value = 69.33309253678253
$ws.Range("J17").Value = 10.43690147392828
$ws.Range("L17").Value = 11.93998820417367

$ws.Range("B18").Value = 29.17411710747698
$ws.Range("C18").Value = 10.44487758185212
$ws.Range("D18").Value = 4.666017177845781
$ws.Range("E18").Value = 9.748176110979266
$ws.Range("F18").Value = 69.13063682386496
$ws.Range("J18").Value = 10.43687081210849
$ws.Range("L18").Value = 11.94789272183035

$ws.Range("B19").Value = 29.15925036506695
$ws.Range("C19").Value = 10.42422466585995
$ws.Range("D19").Value = 4.655061731176797
$ws.Range("E19").Value = 9.749093238601269
$ws.Range("F19").Value = 69.06211988003177
$ws.Range("J19").Value = 10.43686841952421
$ws.Range("L19").Value = 11.95061087480808

$ws.Range("B20").Value = 29.22680539483801
$ws.Range("C20").Value = 10.51720314736405
$ws.Range("D20").Value = 4.704172526115094
$ws.Range("E20").Value = 9.744996295572438
$ws.Range("F20").Value = 69.37057630592614
$ws.Range("J20").Value = 10.43691094565267
$ws.Range("L20").Value = 11.93854512221051

$ws.Range("B21").Value = 29.46504601294622
$ws.Range("C21").Value = 10.82987170027356
$ws.Range("D21").Value = 4.865763572822929
$ws.Range("E21").Value = 9.731758963379859
$ws.Range("F21").Value = 70.40850522423648
$ws.Range("J21").Value = 10.43761088409601
$ws.Range("L21").Value = 11.90091661068078

$ws.Range("B22").Value = 29.62893925813399
$ws.Range("C22").Value = 11.03396295911024
$ws.Range("D22").Value = 4.968724881732784
$ws.Range("E22").Value = 9.723499501008769
$ws.Range("F22").Value = 71.08730405201142
$ws.Range("J22").Value = 10.43848504465946
$ws.Range("L22").Value = 11.87849550388944

$ws.Range("B23").Value = 29.54073077834113
$ws.Range("C23").Value = 10.92508763105673
$ws.Range("D23").Value = 4.914021006983321
$ws.Range("E23").Value = 9.727871968597979
$ws.Range("F23").Value = 70.7250202656702
$ws.Range("J23").Value = 10.43798026187744
$ws.Range("L23").Value = 11.89026327884861

$ws.Range("B24").Value = 29.22305270973305
$ws.Range("C24").Value = 10.51209497761969
$ws.Range("D24").Value = 4.701488163814693
$ws.Range("E24").Value = 9.745219290524117
$ws.Range("F24").Value = 69.35362968172252
$ws.Range("J24").Value = 10.43690651856877
$ws.Range("L24").Value = 11.93919677018268

$ws.Range("B25").Value = 28.91684298281357
$ws.Range("C25").Value = 10.06856034083792
$ws.Range("D25").Value = 4.461565550641629
$ws.Range("E25").Value = 9.758523066601517
$ws.Range("F25").Value = 69.72556874423762
$ws.Range("J25").Value = 10.43752743128698
$ws.Range("L25").Value = 12.00128941426046

